$wb = $excel.ActiveWorkbook

$q = [char]34

# --- Update text values in "Precondiciones" sheet ---
$wsPrecondiciones = $wb.Worksheets.Item("Precondiciones")
$wsPrecondiciones.Range("B4").Value = $q + "CPA_Playa1" + $q + " es el nombre de la playa <Playa1>"

# --- Update text values in "Pasos" sheet ---
$wsPasos = $wb.Worksheets.Item("Pasos")
$wsPasos.Range("B3").Value = "Ingreso " + $q + "CPA_Playa1" + $q + " en el campo nombre de playa"

# --- Update selections on each sheet to match target state ---
$wsDatosGenerales = $wb.Worksheets.Item("DatosGenerales")
$wsDatosGenerales.Activate()
$wsDatosGenerales.Range("B5").Select()

$wsPrecondiciones.Activate()
$wsPrecondiciones.Range("C1").Select()

$wsPasos.Activate()
$wsPasos.Range("B9").Select()

$wsControlCambios = $wb.Worksheets.Item("Control de cambios")
$wsControlCambios.Activate()
$wsControlCambios.Range("F7").Select()

# Restore the originally active sheet/tab
$wsDatosGenerales.Activate()
